# Generate Report for Handback
#
# Appends a new handback row for the localization file
# "652cc506-439c-449b-9f8a-167c2c4f533d" to the Overview sheet and the
# per-locale (zh-cn / de-de) detail sheets, mirroring the layout already
# used by the two existing rows.

$wb = $excel.ActiveWorkbook

$fileId  = "652cc506-439c-449b-9f8a-167c2c4f533d"
$mdName  = "$fileId.md"
$xlfHash = "4b87c66145caa302884654fae3a47a516348e552"
$status  = "Handed back: in sync with en-US"
$reason  = "Include"

# ---------------------------------------------------------------------
# Sheet 1 ("Overview") - one summary row per handed-back file
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/$xlfHash/e2e/$mdName",
    [Type]::Missing,
    [Type]::Missing,
    $mdName) | Out-Null
$wsOverview.Range("A4").Style = "HyperLink"

$wsOverview.Range("B4").Value = $status
$wsOverview.Range("C4").Value = $status

# ---------------------------------------------------------------------
# Sheet 2 ("zh-cn") - handback detail row for the zh-cn locale
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhXlfName = "$fileId.$xlfHash.zh-cn.xlf"
$zhHandoffDateTime  = "2016-03-01 03:30:27"
$zhHandbackDateTime = "2016-03-01 03:31:11"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/$xlfHash/e2e/$mdName",
    [Type]::Missing,
    [Type]::Missing,
    $mdName) | Out-Null
$wsZhCn.Range("A4").Style = "HyperLink"

$wsZhCn.Range("B4").Value = $status

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$xlfHash/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/$zhXlfName",
    [Type]::Missing,
    [Type]::Missing,
    $zhXlfName) | Out-Null
$wsZhCn.Range("C4").Style = "HyperLink"

$wsZhCn.Range("D4").Value = $zhHandoffDateTime
$wsZhCn.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/$xlfHash/e2e/$mdName",
    [Type]::Missing,
    [Type]::Missing,
    $mdName) | Out-Null
$wsZhCn.Range("E4").Style = "HyperLink"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$xlfHash/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/$zhXlfName",
    [Type]::Missing,
    [Type]::Missing,
    $zhXlfName) | Out-Null
$wsZhCn.Range("F4").Style = "HyperLink"

$wsZhCn.Range("G4").Value = $zhHandbackDateTime
$wsZhCn.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Range("H4").Value = $reason

# ---------------------------------------------------------------------
# Sheet 3 ("de-de") - handback detail row for the de-de locale
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deXlfName = "$fileId.$xlfHash.de-de.xlf"
$deHandoffDateTime  = "2016-03-01 03:30:39"
$deHandbackDateTime = "2016-03-01 03:31:31"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/$xlfHash/e2e/$mdName",
    [Type]::Missing,
    [Type]::Missing,
    $mdName) | Out-Null
$wsDeDe.Range("A4").Style = "HyperLink"

$wsDeDe.Range("B4").Value = $status

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$xlfHash/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/$deXlfName",
    [Type]::Missing,
    [Type]::Missing,
    $deXlfName) | Out-Null
$wsDeDe.Range("C4").Style = "HyperLink"

$wsDeDe.Range("D4").Value = $deHandoffDateTime
$wsDeDe.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/$xlfHash/e2e/$mdName",
    [Type]::Missing,
    [Type]::Missing,
    $mdName) | Out-Null
$wsDeDe.Range("E4").Style = "HyperLink"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$xlfHash/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/$deXlfName",
    [Type]::Missing,
    [Type]::Missing,
    $deXlfName) | Out-Null
$wsDeDe.Range("F4").Style = "HyperLink"

$wsDeDe.Range("G4").Value = $deHandbackDateTime
$wsDeDe.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Range("H4").Value = $reason

Write-Host "Added handback row for $fileId to Overview, zh-cn and de-de sheets."
